$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "inntekter"

$ws.Range("A3").Value = "lønn"
$ws.Range("B3").Value = 5000

$ws.Range("A4").Value = "salg av div"
$ws.Range("B4").Value = 2300

$ws.Range("A5").Value = "Inntekt totalt"
$ws.Range("B5").Formula = "=SUM(B3:B4)"

$ws.Range("A7").Value = "Utgifter"
$ws.Range("A8").Value = "Husleie"
$ws.Range("A9").Value = "Drivstoff"
$ws.Range("A10").Value = "Mat"
$ws.Range("A11").Value = "Strøm"
$ws.Range("A12").Value = "Forsikringer"
$ws.Range("A13").Value = "Abonnementer"

$ws.Columns.Item(1).ColumnWidth = 14.166666666666666

$ws.Range("A14").Select()
